$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Time" column (F) entirely, shifting remaining cells left
$ws.Range("F1:F12").Delete()

# Remove rows 7-12 (the trailing combinations that are no longer reported)
$ws.Range("A7:A12").EntireRow.Delete()

# Update the remaining 5 data rows with the new combinations / metrics
$ws.Range("B2").Value = "(7, 'distance', 'manhattan', 'kd_tree', 10)"
$ws.Range("C2").Value = 0.9999095677851279
$ws.Range("D2").Value = 6827.301282051281
$ws.Range("E2").Value = 11.90641025641026

$ws.Range("B3").Value = "(7, 'distance', 'manhattan', 'ball_tree', 30)"
$ws.Range("C3").Value = 0.9999095677851279
$ws.Range("D3").Value = 6827.301282051281
$ws.Range("E3").Value = 11.90641025641026

$ws.Range("B4").Value = "(7, 'distance', 'manhattan', 'ball_tree', 10)"
$ws.Range("C4").Value = 0.9999095677851279
$ws.Range("D4").Value = 6827.301282051281
$ws.Range("E4").Value = 11.90641025641026

$ws.Range("B5").Value = "(7, 'distance', 'manhattan', 'kd_tree', 150)"
$ws.Range("C5").Value = 0.9999095677851279
$ws.Range("D5").Value = 6827.301282051281
$ws.Range("E5").Value = 11.90641025641026

$ws.Range("B6").Value = "(7, 'distance', 'manhattan', 'kd_tree', 30)"
$ws.Range("C6").Value = 0.9999095677851279
$ws.Range("D6").Value = 6827.301282051281
$ws.Range("E6").Value = 11.90641025641026
